$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (default_count, default_value) before the
# existing "most_frequent_value" column (Q), shifting Q:T -> S:V.
$ws.Columns("Q:R").Insert()

# New header cells
$ws.Range("Q1").Value = "default_count"
$ws.Range("R1").Value = "default_value"

# Row 2 - Family Name
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = "<Unspecified>"
$ws.Range("S2").Value = "n/a"

# Row 3
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = "<Unspecified>"
$ws.Range("S3").Value = "n/a"

# Row 4 - Birthplace
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = "<Unspecified>"
$ws.Range("S4").Value = "China"

# Row 5 - Naturalized
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = "<Unspecified>"
$ws.Range("S5").Value = "No"

# Row 6 - Age (type changed from string to int, plus stats recomputed)
$ws.Range("C6").Value = "int"
$ws.Range("E6").Value = 3311
$ws.Range("F6").Value = 147322
$ws.Range("G6").Value = 44.49471458773785
$ws.Range("H6").Value = 45
$ws.Range("I6").Value = 14.46911135959427
$ws.Range("J6").Value = 209.35518353634
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 88
$ws.Range("Q6").Value = 18
$ws.Range("R6").Value = "n/a"
$ws.Range("S6").Value = "50"

# Row 7 - Occupation
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = "<Unspecified>"
$ws.Range("S7").Value = "Labourer"

# Row 8 - Year of Arrival in Queensland
$ws.Range("Q8").Value = 312
$ws.Range("R8").Value = "n/a"
$ws.Range("S8").Value = "n/a"

# Row 9 - Residence
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = "<Unspecified>"
$ws.Range("S9").Value = "Mackay"

# Row 10 - Police District
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = "<Unspecified>"
$ws.Range("S10").Value = "Townsville"

# Row 11 - Police Subdistrict
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = "<Unspecified>"
$ws.Range("S11").Value = "n/a"
